$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (new last row) introduced first -> "Sửa chữa..." becomes shared string index 16
$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "Sửa chữa, viết báo cáo, hoàn thiện"
$ws.Range("C12").Value = "x"
$ws.Range("D12").Value = "x"
$ws.Range("A12").Style = $ws.Range("A9").Style

# Row 5 update -> "Phân tích..." becomes shared string index 17
$ws.Range("B5").Value = "Phân tích các chức năng của từ điển"
$ws.Range("D5").Value = "x"

# Row 6 update -> "Tìm cách..." becomes shared string index 18
$ws.Range("B6").Value = "Tìm cách áp dụng đưa cây nhị phân vào từ điển"
$ws.Range("D6").Value = "x"

# Row 7 update (shift old row5 content)
$ws.Range("B7").Value = "Cài đặt cây nhị phân tìm kiếm"
$ws.Range("C7").Value = "x"
$ws.Range("D7").Value = ""

# Row 8 update (shift old row6 content)
$ws.Range("B8").Value = "Cài đặt từ điển thông thường (tra từ)"
$ws.Range("C8").Value = "x"
$ws.Range("D8").Value = ""

# Row 9 update (shift old row7 content)
$ws.Range("A9").Value = 6
$ws.Range("B9").Value = "Cài đặt tính năng bổ sung từ chưa có trong từ điển"
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = "x"

# Row 10 new (shift old row8 content)
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Cài đặt tính năng cập nhật từ bị sai nghĩa trong từ điển"
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = "x"
$ws.Range("A10").Style = $ws.Range("A9").Style

# Row 11 new (shift old row9 content)
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "Hoàn thiện giao diện "
$ws.Range("C11").Value = "x"
$ws.Range("A11").Style = $ws.Range("A9").Style

# Column G time ranges
$ws.Range("G5").Value = "20/9          -        1/10"
$ws.Range("G6").Value = "1/10          -        20/10"
$ws.Range("G7").Value = "20/9          -        20/10"
$ws.Range("G8").Value = "20/10        -        10/11"
$ws.Range("G9").Value = "10/11        -        30/11"
$ws.Range("G10").Value = "10/11        -        30/11"
$ws.Range("G11").Value = "1/12          -         Cuối khóa"
$ws.Range("G12").Value = "15/12        -         Cuối khóa"

# Column widths / view
$ws.Columns.Item(7).ColumnWidth = 31.15
$ws.Application.ActiveWindow.Zoom = 112
$ws.Range("G12").Select()
